# Update the "Percent Matching" values (column A) on sheet "3B" to reflect
# the corrected/rerun figures from the commit "fix parentheses, rerun files/figures".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3B")

$values = @{
    2  = 49.49832775919732
    3  = 10.7843137254902
    4  = 9.803921568627452
    5  = 11.43790849673203
    6  = 10.7843137254902
    7  = 11.16838487972509
    8  = 9.793814432989691
    9  = 14.26116838487972
    10 = 11.68384879725086
    12 = 8.196721311475411
    13 = 10.81967213114754
    14 = 11.47540983606557
    15 = 10.81967213114754
    16 = 28.57142857142857
    17 = 42.64214046822742
    18 = 14.37908496732026
    19 = 13.0718954248366
    20 = 13.72549019607843
    21 = 14.37908496732026
    22 = 11.51202749140894
    23 = 15.63573883161512
    24 = 13.05841924398625
    25 = 12.54295532646048
    27 = 11.47540983606558
    28 = 9.18032786885246
    29 = 13.44262295081967
    30 = 11.47540983606557
    31 = 39.15966386554622
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
